$d = $word.ActiveDocument

# --- 1. Paragraph 1: add trailing spaces + three red runs ---
$p1 = $d.Paragraphs(1).Range
$p1.InsertAfter("  ")

$endPos = $p1.End - 1
$r2 = $d.Range($endPos, $endPos)
$r2.InsertAfter("(This is a change – Ve")
$r2.Font.Color = 255

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter("rsion for main branch")
$r3.Font.Color = 255

$r4 = $d.Range($r3.End, $r3.End)
$r4.InsertAfter(")")
$r4.Font.Color = 255

# --- 2. Append a new, empty, shaded paragraph at the very end of the doc ---
$endOfDoc = $d.Content.End
$tail = $d.Range($endOfDoc, $endOfDoc)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$tail.InsertXML($xml)
